$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed publish date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$ws.Range("B9").Value = "Alvearie Team"

# The old sheet had a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11). Replace row 10 with a new Jurisdiction row, then remove the
# now-redundant duplicate row 11 entirely (rows below shift up by one).
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Rows.Item(11).Delete()

# Case Sensitive value was blank, now populated ("true" ends up on row 14
# after the row-11 deletion above).
$ws.Range("B14").Value = "true"
